$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("methodNumberOfLines")

$data = @(
  ,@("com.macro.mall.search.domain.EsProductAttributeValue", "EsProductAttributeValue()", "0")
  ,@("com.macro.mall.search.domain.EsProductAttributeValue", "getId()", "3")
  ,@("com.macro.mall.search.domain.EsProductAttributeValue", "setId(java.lang.Long)", "3")
  ,@("com.macro.mall.search.domain.EsProductAttributeValue", "getProductAttributeId()", "3")
  ,@("com.macro.mall.search.domain.EsProductAttributeValue", "setProductAttributeId(java.lang.Long)", "3")
  ,@("com.macro.mall.search.domain.EsProductAttributeValue", "getValue()", "3")
  ,@("com.macro.mall.search.domain.EsProductAttributeValue", "setValue(java.lang.String)", "3")
  ,@("com.macro.mall.search.domain.EsProductAttributeValue", "getType()", "3")
  ,@("com.macro.mall.search.domain.EsProductAttributeValue", "setType(java.lang.Integer)", "3")
  ,@("com.macro.mall.search.domain.EsProductAttributeValue", "getName()", "3")
  ,@("com.macro.mall.search.domain.EsProductAttributeValue", "setName(java.lang.String)", "3")
  ,@("com.macro.mall.search.controller.EsProductController", "EsProductController()", "0")
  ,@("com.macro.mall.search.controller.EsProductController", "importAllList()", "4")
  ,@("com.macro.mall.search.controller.EsProductController", "delete(java.lang.Long)", "4")
  ,@("com.macro.mall.search.controller.EsProductController", "delete(java.util.List)", "4")
  ,@("com.macro.mall.search.controller.EsProductController", "create(java.lang.Long)", "7")
  ,@("com.macro.mall.search.controller.EsProductController", "search(java.lang.String, java.lang.Integer, java.lang.Integer)", "4")
  ,@("com.macro.mall.search.controller.EsProductController", "search(java.lang.String, java.lang.Long, java.lang.Long, java.lang.Integer, java.lang.Integer, java.lang.Integer)", "4")
  ,@("com.macro.mall.search.controller.EsProductController", "recommend(java.lang.Long, java.lang.Integer, java.lang.Integer)", "4")
  ,@("com.macro.mall.search.controller.EsProductController", "searchRelatedInfo(java.lang.String)", "4")
  ,@("com.macro.mall.search.service.EsProductService", "importAll()", "1")
  ,@("com.macro.mall.search.service.EsProductService", "delete(java.lang.Long)", "1")
  ,@("com.macro.mall.search.service.EsProductService", "create(java.lang.Long)", "1")
  ,@("com.macro.mall.search.service.EsProductService", "delete(java.util.List)", "1")
  ,@("com.macro.mall.search.service.EsProductService", "search(java.lang.String, java.lang.Integer, java.lang.Integer)", "1")
  ,@("com.macro.mall.search.service.EsProductService", "search(java.lang.String, java.lang.Long, java.lang.Long, java.lang.Integer, java.lang.Integer, java.lang.Integer)", "1")
  ,@("com.macro.mall.search.service.EsProductService", "recommend(java.lang.Long, java.lang.Integer, java.lang.Integer)", "1")
  ,@("com.macro.mall.search.service.EsProductService", "searchRelatedInfo(java.lang.String)", "1")
  ,@("com.macro.mall.search.config.SwaggerConfig", "SwaggerConfig()", "0")
  ,@("com.macro.mall.search.config.SwaggerConfig", "swaggerProperties()", "3")
  ,@("com.macro.mall.search.config.SwaggerConfig", "springfoxHandlerProviderBeanPostProcessor()", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "EsProduct()", "0")
  ,@("com.macro.mall.search.domain.EsProduct", "getId()", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "setId(java.lang.Long)", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "getProductSn()", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "setProductSn(java.lang.String)", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "getBrandId()", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "setBrandId(java.lang.Long)", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "getBrandName()", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "setBrandName(java.lang.String)", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "getProductCategoryId()", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "setProductCategoryId(java.lang.Long)", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "getProductCategoryName()", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "setProductCategoryName(java.lang.String)", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "getPic()", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "setPic(java.lang.String)", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "getName()", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "setName(java.lang.String)", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "getSubTitle()", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "setSubTitle(java.lang.String)", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "getPrice()", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "setPrice(java.math.BigDecimal)", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "getSale()", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "setSale(java.lang.Integer)", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "getNewStatus()", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "setNewStatus(java.lang.Integer)", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "getRecommandStatus()", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "setRecommandStatus(java.lang.Integer)", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "getStock()", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "setStock(java.lang.Integer)", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "getPromotionType()", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "setPromotionType(java.lang.Integer)", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "getSort()", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "setSort(java.lang.Integer)", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "getAttrValueList()", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "setAttrValueList(java.util.List)", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "getKeywords()", "3")
  ,@("com.macro.mall.search.domain.EsProduct", "setKeywords(java.lang.String)", "3")
  ,@("com.macro.mall.search.service.impl.EsProductServiceImpl", "EsProductServiceImpl()", "0")
  ,@("com.macro.mall.search.service.impl.EsProductServiceImpl", "importAll()", "11")
  ,@("com.macro.mall.search.service.impl.EsProductServiceImpl", "delete(java.lang.Long)", "3")
  ,@("com.macro.mall.search.service.impl.EsProductServiceImpl", "create(java.lang.Long)", "9")
  ,@("com.macro.mall.search.service.impl.EsProductServiceImpl", "delete(java.util.List)", "11")
  ,@("com.macro.mall.search.service.impl.EsProductServiceImpl", "search(java.lang.String, java.lang.Integer, java.lang.Integer)", "4")
  ,@("com.macro.mall.search.service.impl.EsProductServiceImpl", "search(java.lang.String, java.lang.Long, java.lang.Long, java.lang.Integer, java.lang.Integer, java.lang.Integer)", "52")
  ,@("com.macro.mall.search.service.impl.EsProductServiceImpl", "recommend(java.lang.Long, java.lang.Integer, java.lang.Integer)", "34")
  ,@("com.macro.mall.search.service.impl.EsProductServiceImpl", "searchRelatedInfo(java.lang.String)", "16")
  ,@("com.macro.mall.search.service.impl.EsProductServiceImpl", "convertProductRelatedInfo(org.springframework.data.elasticsearch.core.SearchHits)", "37")
  ,@("com.macro.mall.search.domain.EsProductRelatedInfo$ProductAttr", "getAttrId()", "3")
  ,@("com.macro.mall.search.domain.EsProductRelatedInfo$ProductAttr", "setAttrId(java.lang.Long)", "3")
  ,@("com.macro.mall.search.domain.EsProductRelatedInfo$ProductAttr", "getAttrValues()", "3")
  ,@("com.macro.mall.search.domain.EsProductRelatedInfo$ProductAttr", "setAttrValues(java.util.List)", "3")
  ,@("com.macro.mall.search.domain.EsProductRelatedInfo$ProductAttr", "getAttrName()", "3")
  ,@("com.macro.mall.search.domain.EsProductRelatedInfo$ProductAttr", "setAttrName(java.lang.String)", "3")
  ,@("com.macro.mall.search.domain.EsProductRelatedInfo", "EsProductRelatedInfo()", "0")
  ,@("com.macro.mall.search.domain.EsProductRelatedInfo", "getBrandNames()", "3")
  ,@("com.macro.mall.search.domain.EsProductRelatedInfo", "setBrandNames(java.util.List)", "3")
  ,@("com.macro.mall.search.domain.EsProductRelatedInfo", "getProductCategoryNames()", "3")
  ,@("com.macro.mall.search.domain.EsProductRelatedInfo", "setProductCategoryNames(java.util.List)", "3")
  ,@("com.macro.mall.search.domain.EsProductRelatedInfo", "getProductAttrs()", "3")
  ,@("com.macro.mall.search.domain.EsProductRelatedInfo", "setProductAttrs(java.util.List)", "3")
  ,@("com.macro.mall.search.dao.EsProductDao", "getAllEsProductList(java.lang.Long)", "1")
  ,@("com.macro.mall.search.repository.EsProductRepository", "findByNameOrSubTitleOrKeywords(java.lang.String, java.lang.String, java.lang.String, org.springframework.data.domain.Pageable)", "1")
  ,@("com.macro.mall.search.MallSearchApplication", "MallSearchApplication()", "0")
  ,@("com.macro.mall.search.MallSearchApplication", "main(java.lang.String[])", "3")
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $i + 2
  $ws.Cells.Item($r, 1).Value = $data[$i][0]
  $ws.Cells.Item($r, 2).Value = $data[$i][1]
  $cC = $ws.Cells.Item($r, 3)
  $cC.NumberFormat = "@"
  $cC.Value = $data[$i][2]
  $cC.ClearFormats()
}

$wsClassLines = $wb.Worksheets.Item("classNumberOfLines")
$cB13 = $wsClassLines.Cells.Item(13, 2)
$cB13.NumberFormat = "@"
$cB13.Value = "0"
$cB13.ClearFormats()
